# Edit script for LOB1252.xlsx
#
# Fills in the Portuguese "Objetivos:" text on row 10, inserts a dedicated
# "Docentes responsaveis:" data row (professor name moves from row 18 to the
# new row 13), fills the Portuguese "Programa resumido:"/"Programa:" texts,
# restores the "Short syllabus:"/"Syllabus:" rows pushed down by one position,
# reorders "Metodo:" / "Criterio:" / "Norma de recuperacao:" and appends a new
# "Bibliografia:" row with the reference list. Row heights are adjusted to
# match the new content, and rows that no longer need a tall custom height
# are auto-fit back to the default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:) : add the Portuguese objectives paragraph ---
$objetivos = 'Apresentar aos alunos o estudo de tecnologias voltadas à captura de CO2 por meios físicos, químicos e biológicos e comparar os estudos de tecnologias voltadas para à captura de CO2, enfatizando suas vantagens, desvantagens e aplicabilidade quanto ao ponto de vista tecnológico, ambiental e de sustentabilidade.'
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Row 13 : drop "Programa resumido: / Semestral", keep only the professor ---
# --- line that used to be on row 18, and shrink back to the default height ---
$ws.Range("A13").Value = $null
$docente = '5840692 - Diovana Aparecida dos Santos Napoleão'
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
$ws.Rows(13).EntireRow.AutoFit()

# --- Row 14 : "Programa resumido:" gets the new Portuguese summary text ---
$ws.Range("A14").Value = 'Programa resumido:'
$resumo = 'O Ciclo do carbono e emissões. Fontes renováveis e não renováveis. Estudo das emissões na geração de energia. Tecnologias de captura/sequestro e armazenamento de CO2. Alternativas de mitigação de CO2'
$ws.Range("B14").Value = $resumo
$ws.Range("C14").Value = $resumo

# --- Row 15 : "Short syllabus:" (previously row 14), now a 60pt row ---
$ws.Range("A15").Value = 'Short syllabus:'
$shortSyllabus = 'The carbon cycle and emissions. Renewable and non-renewable sources. Study of emissions in energy generation. Capture technologies and storage of CO2. Mitigation alternatives CO2.'
$ws.Range("B15").Value = $shortSyllabus
$ws.Range("C15").Value = $shortSyllabus
$ws.Rows(15).RowHeight = 60

# --- Row 16 : "Programa:" gets the new detailed Portuguese program text ---
$ws.Range("A16").Value = 'Programa:'
$programa = 'Estudo do dióxido de carbono e sua importância na atmosfera. Formação do Ciclo do carbono. Influência do CO2 sobre o meio ambiente. Aquecimento global e os gases de efeito estufa (GEE). Emissões mundiais de GEE por atividades antropogênicas. Emissões históricas globais de CO2 atribuíveis à mudança de uso da Terra. Emissões de CO2 na mudança de uso da Terra em Biomas do Brasil, Potencial de sequestro de carbono em atividades de manejo no Brasil. Fontes renováveis de energia. Termelétricas, Hidrelétricas, Energia eólica, Energia solar, Energia geotérmica, Energia mareomotriz. Fontes não renováveis de energia. O carvão no mundo e no Brasil. Petróleo. Gás natural. Estudo das emissões de CO2 na geração de energia: perspectivas do gerenciamento ambiental para o problema. A questão das emissões de CO2 e a comunidade internacional. Captura, separação e armazenamento de CO2. Estudo das tecnologias avançadas na geração energética para redução das emissões de CO2. Alternativas tecnológicas para a redução das emissões de CO2.'
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Row 17 : "Syllabus:" (previously row 16), stays a 120pt row ---
$ws.Range("A17").Value = 'Syllabus:'
$syllabus = 'Study carbon dioxide and its importance in the atmosphere. Formation of the carbon cycle. Effect of CO2 on the environment. Global warming and greenhouse gases (GHG). GHG global emissions from anthropogenic activities. Global historical emissions CO2 attributable the change in the use Earth in Brazil biomes. Potential sequestration Carbon management activities in Brazil. Renewable energy sources. Thermoelectric, Hydroelectric, Wind Energy, Solar Energy, Geothermal Energy, Energy mareomotriz. Non-renewable energy sources. The coal in the world and in Brazil. Oil. Natural gas. Study of CO2 emissions in energy generation: perspectives of environmental management for the problem. The question of CO2 emissions and the international community. Capture, separation and storage of CO2. Study of advanced technologies in energy generation to reduce CO2 emissions. Technological alternatives to reduce CO2 emissions.'
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus
$ws.Rows(17).RowHeight = 120
# B17 is a brand-new cell (row 17 previously had no column-B content), so make sure
# it picks up the same "normal weight, wrapped text" look used by every other column-B cell.
$ws.Range("B17").Font.Bold = $false
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4160

# --- Row 18 : now just "Avaliacao:" ; the professor pair moved to row 13 ---
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18").Value = $null
$ws.Range("C18").Value = $null
$ws.Rows(18).EntireRow.AutoFit()

# --- Row 19 : "Metodo:" (previously "Criterio:") ---
$ws.Range("A19").Value = 'Método:'

# --- Row 20 : "Criterio:" (previously "Norma de recuperacao:") ---
$ws.Range("A20").Value = 'Critério:'

# --- Row 21 : "Norma de recuperacao:" shrinks from a 120pt row to 60pt ---
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Rows(21).RowHeight = 60

# --- Row 22 (new) : "Bibliografia:" with the reference list, a 120pt row ---
$ws.Range("A22").Value = 'Bibliografia:'
$bibliografia = '1 - Villela, A. A., Freitas, M. A., Rosa, L. P. Emissões de carbono na mudança de uso do solo. Edta Interciência, Vol. 2, 2012.
2 - Marengo, J.A. Mudanças climáticas globais e seus efeitos sobre a biodiversidade. Ministério do Meio Ambiente, 2006.
3 - Patusco, J. A. M. Energia &Desenvolvimento  Ranking Estadual de Energia, Emissões CO2 e socioeconômica. Edta Kiron, 2012.
4 - Goldemberg, J., Palleta, F. C. Energias Renováveis. Edta. Blucher, 2012.
Burattini, M. P. T. C., Energia uma abordagem multidisciplinar. Edta LTF, 2009.
5 - Meyer, L., Pachauri, R. K. Climate Change 2014  Synthesis Report. Intergovernmental Panel on Climate Change. Geneva, Switzerland, 151 pp., 2014.'
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
$ws.Rows(22).RowHeight = 120
# Row 22 is brand new, so B22 needs the usual column-B look explicitly
# (A22 correctly inherits the bold label style, C22 correctly inherits the red style).
$ws.Range("B22").Font.Bold = $false
$ws.Range("B22").WrapText = $true
$ws.Range("B22").VerticalAlignment = -4160

